# Commit: "new results with http1.1 and http2"
#
# For all three result sheets ("Transfer Time (s)", "Throughput (bps)",
# "Overhead Ratio") the old run had three protocol columns:
#   B:C = HTTP/1.1 (Mean/Std Dev)
#   D:E = HTTP/2 SSL (Mean/Std Dev)
#   F:G = HTTP/2 (Mean/Std Dev)
#
# The new run drops the separate "HTTP/2 SSL" vs "HTTP/2" split: the old
# "HTTP/2" columns (F:G) are removed entirely, and the remaining "HTTP/2 SSL"
# header (D1, merged D1:E1) is renamed to plain "HTTP/2". Fresh benchmark
# numbers are written into B4:E7 (sheet 3 only refreshes its D column - the
# HTTP/1.1 and HTTP/2-SSL-derived overhead-ratio values are unchanged).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Transfer Time (s)"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Transfer Time (s)")

$ws.Range("F1:G1").EntireColumn.Delete()
$ws.Range("D1").Value = "HTTP/2"

$ws.Range("B4").Value = 0.01056118428707123
$ws.Range("C4").Value = 0.005242367281925539
$ws.Range("D4").Value = 0.01389548802375793
$ws.Range("E4").Value = 0.01750155598801064

$ws.Range("B5").Value = 0.01183918118476868
$ws.Range("C5").Value = 0.006374975155080995
$ws.Range("D5").Value = 0.04640535354614258
$ws.Range("E5").Value = 0.01160893773458314

$ws.Range("B6").Value = 0.03155416250228882
$ws.Range("C6").Value = 0.01692662311134954
$ws.Range("D6").Value = 0.09788004159927369
$ws.Range("E6").Value = 0.01815001640077183

$ws.Range("B7").Value = 0.3036929368972778
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 3.874284625053406
$ws.Range("E7").Value = 0

# ---------------------------------------------------------------------
# Sheet 2: "Throughput (bps)"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Throughput (bps)")

$ws.Range("F1:G1").EntireColumn.Delete()
$ws.Range("D1").Value = "HTTP/2"

$ws.Range("B4").Value = 9255521.730111934
$ws.Range("C4").Value = 3680421.788750582
$ws.Range("D4").Value = 20744932.44065881
$ws.Range("E4").Value = 14488310.68536471

$ws.Range("B5").Value = 84427327.22572178
$ws.Range("C5").Value = 33649506.05929364
$ws.Range("D5").Value = 23825560.63408853
$ws.Range("E5").Value = 28432632.23320282

$ws.Range("B6").Value = 308957763.0495739
$ws.Range("C6").Value = 95561058.99782108
$ws.Range("D6").Value = 88122599.94104275
$ws.Range("E6").Value = 13777981.31061908

$ws.Range("B7").Value = 276513431.7404806
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 21962448.46551479
$ws.Range("E7").Value = 0

# ---------------------------------------------------------------------
# Sheet 3: "Overhead Ratio"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overhead Ratio")

$ws.Range("F1:G1").EntireColumn.Delete()
$ws.Range("D1").Value = "HTTP/2"

$ws.Range("D4").Value = 1.00908203125
$ws.Range("D5").Value = 1.00091796875
$ws.Range("D6").Value = 1.000090599060059
$ws.Range("D7").Value = 1.000009155273438

Write-Output "edit complete"
